$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (AES) values
$ws.Range("C2").Value = 190000000.0000004
$ws.Range("D2").Value = 105.1172485027671
$ws.Range("E2").Value = 19972277215.52578

# Insert two new rows before the current TOTAL row (row 3) to make room
$ws.Range("A3:A4").EntireRow.Insert()

# New row 3: ISAGEN offer
$ws.Range("A3").Value = "OFERTA"
$ws.Range("B3").Value = "OP1_Wide -ISAGEN"
$ws.Range("C3").Value = 190000000.0000004
$ws.Range("D3").Value = 54.47866916837214
$ws.Range("E3").Value = 10350947141.99072

# New row 4: EPM offer
$ws.Range("A4").Value = "OFERTA"
$ws.Range("B4").Value = "OP1_Wide- EPM"
$ws.Range("C4").Value = 916608868
$ws.Range("D4").Value = 57.83386399818058
$ws.Range("E4").Value = 53011032611.43826

# Row 5 is now the old TOTAL row (shifted down from row 3), update its values
$ws.Range("A5").Value = "TOTAL"
$ws.Range("B5").Value = "TODAS LAS OFERTAS"
$ws.Range("C5").Value = 1296608868.000001
$ws.Range("D5").Value = 64.27092936476407
$ws.Range("E5").Value = 83334256968.95476
